$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Cells.Item(4, 8).Value = 1418.2222
$ws.Cells.Item(4, 9).Value = 1470.5
$ws.Cells.Item(4, 10).Value = 1000
$ws.Cells.Item(4, 11).Value = 1470.5
$ws.Cells.Item(4, 12).Value = 1000
$ws.Cells.Item(4, 13).Value = -1356.5
$ws.Cells.Item(4, 14).Value = -1228
# Row 80
$ws.Cells.Item(80, 8).Value = 48708.215
$ws.Cells.Item(80, 9).Value = 75364.37
$ws.Cells.Item(80, 10).Value = 727.13336
$ws.Cells.Item(80, 11).Value = 226093.11
$ws.Cells.Item(80, 12).Value = 2181.40008
$ws.Cells.Item(80, 13).Value = -225095.11
$ws.Cells.Item(80, 14).Value = -4177.40008
# Row 83
$ws.Cells.Item(83, 8).Value = 48708.215
$ws.Cells.Item(83, 9).Value = 75364.37
$ws.Cells.Item(83, 10).Value = 727.13336
$ws.Cells.Item(83, 11).Value = 678279.33
$ws.Cells.Item(83, 12).Value = 6544.20024
$ws.Cells.Item(83, 13).Value = -673287.33
$ws.Cells.Item(83, 14).Value = -16528.20024
# Row 113
$ws.Cells.Item(113, 8).Value = 9575.182000000001
$ws.Cells.Item(113, 9).Value = 7632.5713
$ws.Cells.Item(113, 11).Value = 7632.5713
$ws.Cells.Item(113, 13).Value = -4378.5713
# Row 138
$ws.Cells.Item(138, 8).Value = 378287.7
$ws.Cells.Item(138, 9).Value = 627830.3
$ws.Cells.Item(138, 10).Value = 3973.7856
$ws.Cells.Item(138, 11).Value = 1883490.9
$ws.Cells.Item(138, 12).Value = 11921.3568
$ws.Cells.Item(138, 13).Value = -1878350.9
$ws.Cells.Item(138, 14).Value = -22201.3568
# Row 141
$ws.Cells.Item(141, 8).Value = 5686.4
$ws.Cells.Item(141, 9).Value = 5633.9546
$ws.Cells.Item(141, 10).Value = 5830.625
$ws.Cells.Item(141, 11).Value = 16901.8638
$ws.Cells.Item(141, 12).Value = 17491.875
$ws.Cells.Item(141, 13).Value = -11721.8638
$ws.Cells.Item(141, 14).Value = -27851.875

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 3706.95
$ws.Cells.Item(2, 9).Value = 3631.7058
$ws.Cells.Item(2, 11).Value = 3631.7058
$ws.Cells.Item(2, 13).Value = -3518.7058
# Row 74
$ws.Cells.Item(74, 8).Value = 6024.0386
$ws.Cells.Item(74, 9).Value = 6616.25
$ws.Cells.Item(74, 11).Value = 6616.25
$ws.Cells.Item(74, 13).Value = -5742.25
# Row 77
$ws.Cells.Item(77, 8).Value = 6024.0386
$ws.Cells.Item(77, 9).Value = 6616.25
$ws.Cells.Item(77, 11).Value = 33081.25
$ws.Cells.Item(77, 13).Value = -28713.25
# Row 116
$ws.Cells.Item(116, 8).Value = 3706.95
$ws.Cells.Item(116, 9).Value = 3631.7058
$ws.Cells.Item(116, 11).Value = 3631.7058
$ws.Cells.Item(116, 13).Value = -1337.7058

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 3706.95
$ws.Cells.Item(3, 9).Value = 3631.7058
$ws.Cells.Item(3, 11).Value = 3631.7058
$ws.Cells.Item(3, 13).Value = -3517.7058
# Row 80
$ws.Cells.Item(80, 8).Value = 1522.8
$ws.Cells.Item(80, 9).Value = 997.5
$ws.Cells.Item(80, 10).Value = 1873
$ws.Cells.Item(80, 11).Value = 997.5
$ws.Cells.Item(80, 12).Value = 1873
$ws.Cells.Item(80, 13).Value = 0.5
$ws.Cells.Item(80, 14).Value = -3869
# Row 83
$ws.Cells.Item(83, 8).Value = 1522.8
$ws.Cells.Item(83, 9).Value = 997.5
$ws.Cells.Item(83, 10).Value = 1873
$ws.Cells.Item(83, 11).Value = 4987.5
$ws.Cells.Item(83, 12).Value = 9365
$ws.Cells.Item(83, 13).Value = 4.5
$ws.Cells.Item(83, 14).Value = -19349
# Row 132
$ws.Cells.Item(132, 8).Value = 78000
$ws.Cells.Item(132, 10).Value = 78000
$ws.Cells.Item(132, 12).Value = 78000
$ws.Cells.Item(132, 14).Value = -88120

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 11274.765
$ws.Cells.Item(31, 9).Value = 15075.667
$ws.Cells.Item(31, 11).Value = 15075.667
$ws.Cells.Item(31, 13).Value = -14780.667
# Row 34
$ws.Cells.Item(34, 8).Value = 11274.765
$ws.Cells.Item(34, 9).Value = 15075.667
$ws.Cells.Item(34, 11).Value = 15075.667
$ws.Cells.Item(34, 13).Value = -14873.667
# Row 62
$ws.Cells.Item(62, 8).Value = 4890.1816
$ws.Cells.Item(62, 9).Value = 2329.6667
$ws.Cells.Item(62, 11).Value = 2329.6667
$ws.Cells.Item(62, 13).Value = -1705.6667
# Row 65
$ws.Cells.Item(65, 8).Value = 4890.1816
$ws.Cells.Item(65, 9).Value = 2329.6667
$ws.Cells.Item(65, 11).Value = 11648.3335
$ws.Cells.Item(65, 13).Value = -8528.333500000001

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Cells.Item(113, 8).Value = 11971.667
$ws.Cells.Item(113, 9).Value = 881.6667
$ws.Cells.Item(113, 11).Value = 2645.0001
$ws.Cells.Item(113, 13).Value = -475.0001000000002

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 6590
$ws.Cells.Item(80, 9).Value = 8387.125
$ws.Cells.Item(80, 10).Value = 2995.75
$ws.Cells.Item(80, 11).Value = 8387.125
$ws.Cells.Item(80, 12).Value = 2995.75
$ws.Cells.Item(80, 13).Value = -7389.125
$ws.Cells.Item(80, 14).Value = -4991.75
# Row 83
$ws.Cells.Item(83, 8).Value = 6590
$ws.Cells.Item(83, 9).Value = 8387.125
$ws.Cells.Item(83, 10).Value = 2995.75
$ws.Cells.Item(83, 11).Value = 41935.625
$ws.Cells.Item(83, 12).Value = 14978.75
$ws.Cells.Item(83, 13).Value = -36943.625
$ws.Cells.Item(83, 14).Value = -24962.75
# Row 102
$ws.Cells.Item(102, 8).Value = 7595.48
$ws.Cells.Item(102, 9).Value = 8737.833000000001
$ws.Cells.Item(102, 11).Value = 8737.833000000001
$ws.Cells.Item(102, 13).Value = -7115.833000000001
# Row 113
$ws.Cells.Item(113, 8).Value = 11487.833
$ws.Cells.Item(113, 9).Value = 29005
$ws.Cells.Item(113, 10).Value = 2729.25
$ws.Cells.Item(113, 11).Value = 29005
$ws.Cells.Item(113, 12).Value = 2729.25
$ws.Cells.Item(113, 13).Value = -26835
$ws.Cells.Item(113, 14).Value = -7069.25
# Row 132
$ws.Cells.Item(132, 8).Value = 6931.778
$ws.Cells.Item(132, 9).Value = 7385.875
$ws.Cells.Item(132, 10).Value = 3299
$ws.Cells.Item(132, 11).Value = 22157.625
$ws.Cells.Item(132, 12).Value = 9897
$ws.Cells.Item(132, 13).Value = -19627.625
$ws.Cells.Item(132, 14).Value = -14957

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 12288.889
$ws.Cells.Item(22, 9).Value = 20480
$ws.Cells.Item(22, 10).Value = 2050
$ws.Cells.Item(22, 11).Value = 20480
$ws.Cells.Item(22, 12).Value = 2050
$ws.Cells.Item(22, 13).Value = -20185
$ws.Cells.Item(22, 14).Value = -2640
# Row 27
$ws.Cells.Item(27, 8).Value = 12288.889
$ws.Cells.Item(27, 9).Value = 20480
$ws.Cells.Item(27, 10).Value = 2050
$ws.Cells.Item(27, 11).Value = 20480
$ws.Cells.Item(27, 12).Value = 2050
$ws.Cells.Item(27, 13).Value = -20373
$ws.Cells.Item(27, 14).Value = -2264
# Row 40
$ws.Cells.Item(40, 8).Value = 21123.576
$ws.Cells.Item(40, 9).Value = 25702.264
$ws.Cells.Item(40, 11).Value = 25702.264
$ws.Cells.Item(40, 13).Value = -25566.264
# Row 122
$ws.Cells.Item(122, 8).Value = 5781.607
$ws.Cells.Item(122, 9).Value = 5680.9443
$ws.Cells.Item(122, 11).Value = 17042.8329
$ws.Cells.Item(122, 13).Value = -14592.8329
# Row 136
$ws.Cells.Item(136, 8).Value = 5119.2607
$ws.Cells.Item(136, 9).Value = 4367.2354
$ws.Cells.Item(136, 11).Value = 13101.7062
$ws.Cells.Item(136, 13).Value = -10551.7062

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Cells.Item(126, 8).Value = 20336.182
$ws.Cells.Item(126, 9).Value = 23916.889
$ws.Cells.Item(126, 11).Value = 71750.667
$ws.Cells.Item(126, 13).Value = -69280.667
# Row 136
$ws.Cells.Item(136, 8).Value = 392266.25
$ws.Cells.Item(136, 9).Value = 471813.62
$ws.Cells.Item(136, 11).Value = 1415440.86
$ws.Cells.Item(136, 13).Value = -1412890.86
